# Updates the "Pais" (countries COVID stats) sheet with refreshed data
# (commit message: "Update countries & provincias Spain").
#
# Two kinds of changes are applied:
#   1. Simple numeric refreshes for rows whose country didn't move in the
#      table (Estados Unidos, Espana, Alemania, Corea del Sur, Costa Rica,
#      Libia).
#   2. Three "new country row inserted, older rows cascade down one slot"
#      groups (Niger before Ghana, Guinea before Isla de Man, Bahamas
#      before Puerto Rico); the now-stale row that used to hold that
#      country's old data drops off the bottom of each cascade. None of
#      the underlying worksheet rows actually move, so this is applied by
#      rewriting, top-to-bottom, the country name (column A) and stats
#      (columns B-H) for every row touched by each cascade.
#   3. The "data as of" timestamp banner in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Country, $Stats) {
    if ($Country -ne $null) {
        $ws.Cells.Item($Row, 1).Value = $Country
    }
    for ($i = 0; $i -lt $Stats.Length; $i++) {
        $ws.Cells.Item($Row, 2 + $i).Value = $Stats[$i]
    }
}

# --- 1. Plain numeric refreshes (country stays put) -----------------------
Set-Row 4   $null @(425469, 25134, 22202, 388677, 9234, 1749, 14590)  # Estados Unidos
Set-Row 5   $null @(146824, 4882,  48021, 84118,  7069, 640,  14685) # Espana
Set-Row 8   $null @(112113, 4450,  36081, 73824,  4895, 192,  2208)  # Alemania
Set-Row 19  $null @(12942,  303,   4512,  8157,   267,  30,   273)   # Corea del Sur
Set-Row 87  $null @(502,    19,    29,    470,    15,   1,    3)     # Costa Rica
Set-Row 163 $null @(21,     1,     8,     12,     0,    0,    1)     # Libia

# --- 2a. Niger inserted before Ghana; Ghana/Honduras/Malta/San Marino -----
#         cascade down one row; stale Niger row (ex row 102) dropped.
Set-Row 98  "Niger"      @(342, 64, 28, 303, 0,  0, 11)
Set-Row 99  "Ghana"      @(313, 26, 34, 273, 2,  1, 6)
Set-Row 100 "Honduras"   @(312, 7,  6,  284, 10, 0, 22)
Set-Row 101 "Malta"      @(299, 6,  16, 282, 4,  1, 1)
Set-Row 102 "San Marino" @(279, 0,  40, 205, 14, 0, 34)

# --- 2b. Guinea inserted before Isla de Man; Isla de Man/Martinica --------
#         cascade down one row; stale Guinea row (ex row 121) dropped.
Set-Row 119 "Guinea"      @(164, 20, 5,  159, 0,  0, 0)
Set-Row 120 "Isla de Man" @(158, 8,  82, 75,  7,  0, 1)
Set-Row 121 "Martinica"   @(154, 2,  50, 98,  19, 2, 6)

# --- 2c. Bahamas inserted before Puerto Rico; Puerto Rico/Zambia/----------
#         Bermudas/Guyana cascade down one row; stale Bahamas row
#         (ex row 152) dropped.
Set-Row 148 "Bahamas"     @(40, 7, 5,  28, 1, 1, 7)
Set-Row 149 "Puerto Rico" @(39, 0, 1,  36, 0, 0, 2)
Set-Row 150 "Zambia"      @(39, 0, 7,  31, 1, 0, 1)
Set-Row 151 "Bermudas"    @(39, 0, 17, 20, 0, 0, 2)
Set-Row 152 "Guyana"      @(37, 4, 8,  23, 4, 1, 6)

# --- 3. Timestamp banner ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 23:22"
